$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1, matching the formatting of the existing header cells (E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "OSMO_DEF"

# New data cell F2 (plain, unstyled like C2/E2)
$ws.Range("F2").Value = "[]"
